# Sample Bnumber List.xlsx
# Re-populate the header/data table with the new sample names so the
# "first name" matching feature can be tested against real-looking
# first/middle/last name + B-number rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (text itself is unchanged, just re-asserted)
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Middle Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "B Number"

# Row 2
$ws.Range("A2").Value = "Jim "
$ws.Range("B2").Value = "Aiden"
$ws.Range("C2").Value = "Test"
$ws.Range("D2").Value = 1

# Row 3
$ws.Range("A3").Value = "Joe "
$ws.Range("B3").Value = "Jaundice"
$ws.Range("C3").Value = "Test"
$ws.Range("D3").Value = 2

# Row 4
$ws.Range("A4").Value = "Jane "
$ws.Range("B4").Value = "Maleficent"
$ws.Range("C4").Value = "Test"
$ws.Range("D4").Value = 3

# Column A widened after the new (shorter) first names were typed in
$ws.Columns.Item(1).ColumnWidth = 12.5

# Selection left on C11 when the workbook was saved
$ws.Range("C11").Select()
